# Adds the "company / insurance / species / debtor / debt" metadata columns
# to the 保險 (Insurance) and 債務 (Debt) sheets, matching the other sheets'
# shape (property_category, category, date, legislator_name, legislator_id,
# source_file, index [+ owner/total/register_date/register_reason for debt]).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (Insurance) -> add columns E..K
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item("保險")

# Copy existing header/data formatting onto the new columns first so the
# new cells pick up the same style (bold/border header vs plain data).
$ins.Range("B1:D1").Copy()
$ins.Range("E1:K1").PasteSpecial(-4122)
$ins.Range("B2:D7").Copy()
$ins.Range("E2:K7").PasteSpecial(-4122)

# Header row
$ins.Cells.Item(1,5).Value  = "property_category"
$ins.Cells.Item(1,6).Value  = "category"
$ins.Cells.Item(1,7).Value  = "date"
$ins.Cells.Item(1,8).Value  = "legislator_name"
$ins.Cells.Item(1,9).Value  = "legislator_id"
$ins.Cells.Item(1,10).Value = "source_file"
$ins.Cells.Item(1,11).Value = "index"

# Data rows 2-7: same constant metadata on every row, "index" = column A value
$insIndex = @(146,147,148,149,150,151)
for ($i = 0; $i -lt $insIndex.Count; $i++) {
    $r = $i + 2
    $ins.Cells.Item($r,5).Value  = "insurance"
    $ins.Cells.Item($r,6).Value  = "normal"
    $ins.Cells.Item($r,7).Value  = "2011-11-17"
    $ins.Cells.Item($r,8).Value  = "吳育昇"
    $ins.Cells.Item($r,9).Value  = 1322
    $ins.Cells.Item($r,10).Value = "tmpe6fb1"
    $ins.Cells.Item($r,11).Value = $insIndex[$i]
}

# ---------------------------------------------------------------------
# Sheet "債務" (Debt) -> add columns H..N
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")

$debt.Range("B1:G1").Copy()
$debt.Range("H1:N1").PasteSpecial(-4122)
$debt.Range("B2:G2").Copy()
$debt.Range("H2:N2").PasteSpecial(-4122)

# Header row
$debt.Cells.Item(1,8).Value  = "property_category"
$debt.Cells.Item(1,9).Value  = "category"
$debt.Cells.Item(1,10).Value = "date"
$debt.Cells.Item(1,11).Value = "legislator_name"
$debt.Cells.Item(1,12).Value = "legislator_id"
$debt.Cells.Item(1,13).Value = "source_file"
$debt.Cells.Item(1,14).Value = "index"

# Data row 2
$debt.Cells.Item(2,8).Value  = "debt"
$debt.Cells.Item(2,9).Value  = "normal"
$debt.Cells.Item(2,10).Value = "2011-11-17"
$debt.Cells.Item(2,11).Value = "吳育昇"
$debt.Cells.Item(2,12).Value = 1322
$debt.Cells.Item(2,13).Value = "tmpe6fb1"
$debt.Cells.Item(2,14).Value = 161

# The "total" column (E) was stored as a shared-string "1938296" before;
# it should be a genuine number.
$debt.Cells.Item(2,5).Value = 1938296

Write-Output "done"
